$d = $word.ActiveDocument

# Insert the new Lesson 9 content at the very end of the document body,
# right after the last existing paragraph ("    })"). Using a Range
# collapsed at Content.End + InsertXML appends brand-new paragraphs
# without disturbing the pre-existing final paragraph, and (unlike
# InsertParagraphAfter) does not leave a stray empty <w:r/> behind.
$insertionXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Lesson</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> 9 </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">– </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>#</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"># JUGGLING ASYNC </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr><w:r><w:t xml:space="preserve"> This problem is the same as the previous problem (HTTP COLLECT) in that you need to use </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>http.get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>). However, this time you will be provided with three URLs as the first three command-line arguments.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr><w:r><w:t xml:space="preserve"> You must collect the complete content provided to you by each of the URLs and print it to the console (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stdout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>). You don't need to print out the length, just the data as a String; one line per URL. The catch is that you must print them out in the same order as the URLs are provided to you as command-line arguments.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing/></w:pPr></w:p>
'@

$r = $d.Range($d.Content.End, $d.Content.End)
$null = $r.InsertXML($insertionXml)

Write-Output "Paragraph count after insert: $($d.Paragraphs.Count)"
